$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.367401
$ws.Range("H2").Value = 37.102203
$ws.Range("I2").Value = 0.8693283326101076
$ws.Range("J2").Value = 0.8693283326101076
$ws.Range("M2").Value = 48.42420966666666
$ws.Range("N2").Value = 145.272629
$ws.Range("O2").Value = 0.6311762527593259
$ws.Range("P2").Value = 0.6311762527593258
$ws.Range("Q2").Value = 598.881619055743
$ws.Range("R2").Value = 5389.934571501687
$ws.Range("S2").Value = 0.5486993993943605
$ws.Range("T2").Value = 0.5486993993943604

$ws.Range("G3").Value = 12.367401
$ws.Range("H3").Value = 37.102203
$ws.Range("I3").Value = 0.8693283326101076
$ws.Range("J3").Value = 0.8693283326101076
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("O3").Value = 0.08928392431779728
$ws.Range("P3").Value = 0.08928392431779726
$ws.Range("Q3").Value = 84.71564149844801
$ws.Range("R3").Value = 762.4407734860321
$ws.Range("S3").Value = 0.07761704505607774
$ws.Range("T3").Value = 0.07761704505607772

$ws.Range("G4").Value = 12.367401
$ws.Range("H4").Value = 37.102203
$ws.Range("I4").Value = 0.8693283326101076
$ws.Range("J4").Value = 0.8693283326101076
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2795398229228769
$ws.Range("P4").Value = 0.2795398229228769
$ws.Range("Q4").Value = 265.236946115658
$ws.Range("R4").Value = 2387.132515040922
$ws.Range("S4").Value = 0.2430118881596693
$ws.Range("T4").Value = 0.2430118881596693

$ws.Range("H5").Value = 4.303227
$ws.Range("I5").Value = 0.1008273593013545
$ws.Range("J5").Value = 0.1008273593013546
$ws.Range("M5").Value = 48.42420966666666
$ws.Range("N5").Value = 145.272629
$ws.Range("O5").Value = 0.6311762527593259
$ws.Range("P5").Value = 0.6311762527593258
$ws.Range("Q5").Value = 69.46012216375365
$ws.Range("R5").Value = 625.1410994737829
$ws.Range("S5").Value = 0.06363983481944713
$ws.Range("T5").Value = 0.06363983481944711

$ws.Range("H6").Value = 4.303227
$ws.Range("I6").Value = 0.1008273593013545
$ws.Range("J6").Value = 0.1008273593013546
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("Q6").Value = 9.825579247098666
$ws.Range("R6").Value = 88.43021322388799
$ws.Range("S6").Value = 0.009002262317025491
$ws.Range("T6").Value = 0.009002262317025491

$ws.Range("H7").Value = 4.303227
$ws.Range("I7").Value = 0.1008273593013545
$ws.Range("J7").Value = 0.1008273593013546
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2795398229228769
$ws.Range("P7").Value = 0.2795398229228769
$ws.Range("Q7").Value = 30.76299237332199
$ws.Range("R7").Value = 276.8669313598979
$ws.Range("S7").Value = 0.02818526216488193
$ws.Range("T7").Value = 0.02818526216488194

$ws.Range("I8").Value = 0.02984430808853782
$ws.Range("J8").Value = 0.02984430808853782
$ws.Range("M8").Value = 48.42420966666666
$ws.Range("N8").Value = 145.272629
$ws.Range("O8").Value = 0.6311762527593259
$ws.Range("P8").Value = 0.6311762527593258
$ws.Range("Q8").Value = 20.55978952624111
$ws.Range("R8").Value = 185.03810573617
$ws.Range("S8").Value = 0.01883701854551814
$ws.Range("T8").Value = 0.01883701854551814

$ws.Range("I9").Value = 0.02984430808853782
$ws.Range("J9").Value = 0.02984430808853782
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("O9").Value = 0.08928392431779728
$ws.Range("P9").Value = 0.08928392431779726
$ws.Range("S9").Value = 0.002664616944694036
$ws.Range("T9").Value = 0.002664616944694036

$ws.Range("I10").Value = 0.02984430808853782
$ws.Range("J10").Value = 0.02984430808853782
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2795398229228769
$ws.Range("P10").Value = 0.2795398229228769
$ws.Range("Q10").Value = 9.105665649446665
$ws.Range("R10").Value = 81.95099084501999
$ws.Range("S10").Value = 0.008342672598325645
$ws.Range("T10").Value = 0.008342672598325645

